$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The roster was reshuffled/updated: several players swapped rows, two
# players were dropped/replaced ("Zach Edey" -> "Jose Alvarado"), and
# each row's position/team was refreshed to match the (possibly new)
# player in that row.

$ws.Range("A2").Value = "James Harden"
$ws.Range("B2").Value = "PG,SG"
$ws.Range("C2").Value = "LA Clippers"

$ws.Range("A3").Value = "Amen Thompson"
$ws.Range("B3").Value = "SG,SF,PF"
$ws.Range("C3").Value = "Houston Rockets"

$ws.Range("A4").Value = "Keyonte George"
$ws.Range("B4").Value = "PG,SG"
$ws.Range("C4").Value = "Utah Jazz"

$ws.Range("A5").Value = "Jose Alvarado"
$ws.Range("B5").Value = "PG"
$ws.Range("C5").Value = "New Orleans Pelicans"

$ws.Range("A6").Value = "Jayson Tatum"
$ws.Range("B6").Value = "SF,PF"
$ws.Range("C6").Value = "Boston Celtics"

$ws.Range("A7").Value = "Zion Williamson"
$ws.Range("B7").Value = "PF,C"
$ws.Range("C7").Value = "New Orleans Pelicans"

$ws.Range("A8").Value = "Anthony Edwards"
$ws.Range("B8").Value = "SG,SF"
$ws.Range("C8").Value = "Minnesota Timberwolves"

$ws.Range("A9").Value = "RJ Barrett"
$ws.Range("B9").Value = "SG,SF,PF"
$ws.Range("C9").Value = "Toronto Raptors"

$ws.Range("A10").Value = "Ivica Zubac"
$ws.Range("B10").Value = "C"
$ws.Range("C10").Value = "LA Clippers"

$ws.Range("A11").Value = "Draymond Green"
$ws.Range("B11").Value = "PF,C"
$ws.Range("C11").Value = "Golden State Warriors"

$ws.Range("A12").Value = "Bobby Portis"
$ws.Range("B12").Value = "PF,C"
$ws.Range("C12").Value = "Milwaukee Bucks"

$ws.Range("A13").Value = "Nicolas Claxton"
$ws.Range("B13").Value = "C"
$ws.Range("C13").Value = "Brooklyn Nets"

$ws.Range("A14").Value = "Paul George"
$ws.Range("B14").Value = "SG,SF,PF"
$ws.Range("C14").Value = "Philadelphia 76ers"

$ws.Range("A15").Value = "Anfernee Simons"
$ws.Range("B15").Value = "PG,SG"
$ws.Range("C15").Value = "Portland Trail Blazers"

$ws.Range("A16").Value = "Jaren Jackson Jr."
$ws.Range("B16").Value = "PF,C"
$ws.Range("C16").Value = "Memphis Grizzlies"

$ws.Range("A17").Value = "Giannis Antetokounmpo"
$ws.Range("B17").Value = "PF,C"
$ws.Range("C17").Value = "Milwaukee Bucks"

$ws.Range("A18").Value = "Fred VanVleet"
$ws.Range("B18").Value = "PG"
$ws.Range("C18").Value = "Houston Rockets"

$ws.Range("A19").Value = "Bradley Beal"
$ws.Range("B19").Value = "PG,SG,SF"
$ws.Range("C19").Value = "Phoenix Suns"
